$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1 using the same formatting (style) as the
# existing header cell H1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in new data columns I (I0) and J (IF) for rows 2-8
$iValues = @(8, 1, 1, 1, 4, 4, 6)
$jValues = @(8, 5, 4, 5, 5, 5, 6)

for ($i = 0; $i -lt 7; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
